$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.726.70"
$ws.Range("E2").Value = "  -1.41%  "
$ws.Range("D3").Value = "3.268.88"
$ws.Range("E3").Value = "  -0.37%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'581.28"
$ws.Range("E5").Value = "  -0.33%  "
$ws.Range("D6").Value = "'184.22"
$ws.Range("E6").Value = "  +1.22%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "'0.604"
$ws.Range("E8").Value = "  +1.05%  "
$ws.Range("D9").Value = "'0.130"
$ws.Range("E9").Value = "  -3.21%  "
$ws.Range("E10").Value = "  -1.64%  "
$ws.Range("E11").Value = "  -3.03%  "
$ws.Range("D12").Value = "3.833.10"
$ws.Range("E12").Value = "  -0.72%  "
$ws.Range("E13").Value = "  +0.70%  "
$ws.Range("D14").Value = "'27.32"
$ws.Range("E14").Value = "  -5.02%  "
$ws.Range("D15").Value = "67.789.58"
$ws.Range("E15").Value = "  -1.32%  "
$ws.Range("D16").Value = "'0.0000168"
$ws.Range("E16").Value = "  -1.89%  "
$ws.Range("D17").Value = "3.263.14"
$ws.Range("E17").Value = "  -0.35%  "
$ws.Range("D18").Value = "'5.72"
$ws.Range("E18").Value = "  -1.63%  "
$ws.Range("D19").Value = "'13.45"
$ws.Range("E19").Value = "  -0.90%  "
$ws.Range("D20").Value = "'400.84"
$ws.Range("E20").Value = "  +1.82%  "
$ws.Range("D21").Value = "'7.57"
$ws.Range("E21").Value = "  -1.88%  "
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("D23").Value = "'70.95"
$ws.Range("E23").Value = "  -0.97%  "
$ws.Range("E24").Value = "  -1.31%  "
$ws.Range("E25").Value = "  -2.59%  "
$ws.Range("D26").Value = "'0.188"
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("D27").Value = "'9.55"
$ws.Range("E27").Value = "  -0.67%  "
$ws.Range("E28").Value = "  +0.63%  "
$ws.Range("E29").Value = "  -1.68%  "
$ws.Range("D30").Value = "'22.68"
$ws.Range("E30").Value = "  -1.41%  "
$ws.Range("D31").Value = "'5.47"
$ws.Range("E31").Value = "  -4.34%  "
$ws.Range("D32").Value = "'6.94"
$ws.Range("E32").Value = "  -2.81%  "
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("E34").Value = "  -3.62%  "
$ws.Range("D35").Value = "'163.69"
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("D36").Value = "'1.46"
$ws.Range("E36").Value = "  -4.13%  "
$ws.Range("D37").Value = "'1.89"
$ws.Range("E37").Value = "  +1.05%  "
$ws.Range("D38").Value = "'26.99"
$ws.Range("E38").Value = "  +3.39%  "
$ws.Range("E39").Value = "  -3.44%  "
$ws.Range("D40").Value = "'4.51"
$ws.Range("E40").Value = "  -1.71%  "
$ws.Range("D41").Value = "2.678.21"
$ws.Range("E41").Value = "  +2.82%  "
$ws.Range("D42").Value = "'6.29"
$ws.Range("E42").Value = "  -3.78%  "
$ws.Range("D43").Value = "'40.73"
$ws.Range("E43").Value = "  -1.90%  "
$ws.Range("E44").Value = "  -1.14%  "
$ws.Range("D45").Value = "'2.43"
$ws.Range("E45").Value = "  -4.85%  "
$ws.Range("D46").Value = "'336.52"
$ws.Range("E46").Value = "  -1.92%  "
$ws.Range("D47").Value = "'24.54"
$ws.Range("E47").Value = "  -0.90%  "
$ws.Range("D48").Value = "'0.0274"
$ws.Range("E48").Value = "  -2.38%  "
$ws.Range("E49").Value = "  -0.20%  "
$ws.Range("E50").Value = "  -1.35%  "
$ws.Range("D51").Value = "'0.968"
$ws.Range("E51").Value = "  -1.54%  "
